$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.744.67"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "2.207.03"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'229.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "'0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.15%  "
$ws.Range("D7").Value = "'60.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.23%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.402"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").Value = "'57.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("D13").Value = "2.535.48"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "'15.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("D15").Value = "'22.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  -4.06%  "
$ws.Range("D18").Value = "2.214.27"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "41.714.35"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'72.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  -4.20%  "
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "'242.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("E26").Value = "  -3.26%  "
$ws.Range("D27").Value = "'9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "'169.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  -5.81%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'19.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("D32").Value = "'2.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.36%  "
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("D34").Value = "'5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").Value = "'4.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").Value = "'2.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("E38").Value = "  -8.79%  "
$ws.Range("E39").Value = "  -8.97%  "
$ws.Range("D40").Value = "'0.000240"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.42%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").Value = "'8.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("D44").Value = "'0.0954"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("D46").Value = "'97.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("E47").Value = "  -14.69%  "
$ws.Range("D48").Value = "1.467.54"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("D49").Value = "'16.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.53%  "
$ws.Range("D50").Value = "'2.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  -5.62%  "
